$wb = $excel.ActiveWorkbook
$wsOptions = $wb.Worksheets.Item("Geometry options")
$wsCalc    = $wb.Worksheets.Item("geometry calculation")

# --- Fix up H15 on 'geometry calculation' BEFORE changing the lookup key / ---
# --- source data, so the dependent cells (K15, etc.) recalc against the  ---
# --- new formula rather than a stale cached value.                       ---
# H15 used to be a hard-coded "5"; it should become a VLOOKUP like its
# neighbours H12:H14, picking up the matching style (s="45", same as H14).
$wsCalc.Range("H14").Copy()
$wsCalc.Range("H15").PasteSpecial(-4122)  # xlPasteFormats
$wsCalc.Range("H15").Formula = "=VLOOKUP(`$B`$2,'Geometry options'!B4:O39,10,FALSE)"

# --- Data corrections on 'Geometry options' row 10 (the "BC-LEEP" row) ---
$wsOptions.Range("K10").Value = 0.001
$wsOptions.Range("M10").Value = 4

# --- Switch the selected archetype on 'geometry calculation' from ---
# --- "NZEH-Arch" to "BC-LEEP" -- this drives every VLOOKUP-based    ---
# --- geometry figure (and the downstream retrofit-cost sheet) to    ---
# --- recompute against the new row.                                 ---
$wsCalc.Range("B2").Value = "BC-LEEP"

# --- Restore the selections recorded in the saved view state ---
$wsCalc.Activate()
$wsCalc.Range("F16").Select()

$wsOptions.Activate()
$wsOptions.Range("B11").Select()
